{"js": "// Change the year in the astromap link from 2018 to 2022.\n//\n// Before the edit the paragraph reads, spread across three differently\n// formatted runs:\n//   \"(\" + \"http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/\" (Hyperlink\n//   character style) + \").\"\n// all three carrying explicit Calibri/CastleT-Book fonts and size 19.\n//\n// After the edit it is a single plain run (no character formatting, no\n// Hyperlink style) holding the whole string with the year updated, preceded\n// by an empty run:\n//   <w:r/><w:r><w:t>(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).</w:t></w:r>\n\nconst oldUrl = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).\";\nconst newUrl = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\nconst body = context.document.body;\nconst results = body.search(oldUrl, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the astromap link text to update.\");\n}\n\nconst target = results.items[0];\n\n// Replace the whole matched range with fresh, plain-formatted OOXML: a\n// leading empty run followed by one run carrying the updated link text.\n// Using insertOoxml (rather than insertText) lands the replacement without\n// inheriting the old runs' character formatting (fonts/size/Hyperlink\n// style), matching how the text collapsed into a single, unformatted run.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p><w:r/><w:r><w:t>\" +\n  escapeXml(newUrl) +\n  \"</w:t></w:r></w:p></w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Change the year in the astromap link from 2018 to 2022.\n#\n# The paragraph currently reads, spread across three differently-formatted\n# runs: \"(\" + \"http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/\"\n# (Hyperlink character style) + \").\" \u2014 all three runs additionally carry\n# Calibri/CastleT-Book fonts and size 19.\n#\n# We select the whole \"(...).\" text, delete it, and insert the replacement\n# text with the updated year. Deleting the old (variously formatted) runs\n# and inserting fresh text collapses everything into a single run that\n# picks up the paragraph mark's default run formatting instead of the\n# old runs' explicit character formatting / Hyperlink style \u2014 matching\n# how the edit turned three heavily-formatted runs into one plain run.\n\n$d = $word.ActiveDocument\n\n$oldText = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).\"\n$newText = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$found = $find.Execute()\n\nif ($found) {\n    $rng = $find.Parent\n    $rng.Delete()\n    $rng.InsertAfter($newText)\n} else {\n    throw \"Could not find the astromap link text to update.\"\n}\n"}
